$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.211.68"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "1.826.70"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'235.48"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("D6").Value = "'0.6104"
$ws.Range("E6").Value = "  -2.92%  "

$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("D8").Value = "'0.07092"
$ws.Range("E8").Value = "  -4.57%  "

$ws.Range("D9").Value = "'0.2805"
$ws.Range("E9").Value = "  -2.98%  "

$ws.Range("D10").Value = "'23.47"
$ws.Range("E10").Value = "  -5.98%  "

$ws.Range("D11").Value = "'0.07671"
$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("D12").Value = "1.824.18"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("D13").Value = "'4.806"
$ws.Range("E13").Value = "  -3.05%  "

$ws.Range("D14").Value = "'0.000009997"
$ws.Range("E14").Value = "  -1.94%  "

$ws.Range("D15").Value = "'0.6314"
$ws.Range("E15").Value = "  -6.32%  "

$ws.Range("D16").Value = "2.064.38"
$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("D17").Value = "'78.55"
$ws.Range("E17").Value = "  -3.69%  "

$ws.Range("D18").Value = "'5.859"
$ws.Range("E18").Value = "  -5.78%  "

$ws.Range("D19").Value = "29.178.51"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").Value = "'226.81"
$ws.Range("E20").Value = "  -0.92%  "

$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("E22").Value = "  -4.19%  "

$ws.Range("D23").Value = "'6.998"
$ws.Range("E23").Value = "  -4.54%  "

$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").Value = "'155.79"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("D26").Value = "'8.043"
$ws.Range("E26").Value = "  -4.97%  "

$ws.Range("D27").Value = "'0.1306"
$ws.Range("E27").Value = "  -2.87%  "

$ws.Range("D28").Value = "'16.57"
$ws.Range("E28").Value = "  -4.54%  "

$ws.Range("D29").Value = "'1.490"
$ws.Range("E29").Value = "  +1.81%  "

$ws.Range("D30").Value = "'0.06340"
$ws.Range("E30").Value = "  -14.38%  "

$ws.Range("D31").Value = "'1.453"
$ws.Range("E31").Value = "  -1.50%  "

$ws.Range("D32").Value = "'3.823"
$ws.Range("E32").Value = "  -5.26%  "

$ws.Range("D33").Value = "'3.804"
$ws.Range("E33").Value = "  -5.71%  "

$ws.Range("D34").Value = "'1.124"
$ws.Range("E34").Value = "  -1.20%  "

$ws.Range("D35").Value = "'1.741"
$ws.Range("E35").Value = "  -4.32%  "

$ws.Range("D36").Value = "'0.6442"
$ws.Range("E36").Value = "  -7.21%  "

$ws.Range("D37").Value = "'2.546"
$ws.Range("E37").Value = "  -1.27%  "

$ws.Range("D38").Value = "1.215.45"
$ws.Range("E38").Value = "  -1.30%  "

$ws.Range("D39").Value = "'2.727"
$ws.Range("E39").Value = "  -2.61%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01738"
$ws.Range("E40").Value = "  -5.39%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.540"
$ws.Range("E41").Value = "  -4.46%  "

$ws.Range("D42").Value = "'0.9130"
$ws.Range("E42").Value = "  -2.04%  "

$ws.Range("D43").Value = "'1.003"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").Value = "'101.00"
$ws.Range("E44").Value = "  +0.46%  "

$ws.Range("D45").Value = "1.976.62"
$ws.Range("E45").Value = "  -0.51%  "

$ws.Range("D46").Value = "'62.64"
$ws.Range("E46").Value = "  -4.01%  "

$ws.Range("E47").Value = "  -4.08%  "

$ws.Range("D48").Value = "'8.556"
$ws.Range("E48").Value = "  -3.91%  "

$ws.Range("D49").Value = "'1.603"
$ws.Range("E49").Value = "  -5.85%  "

$ws.Range("D50").Value = "'0.4571"
$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").Value = "'0.05521"
$ws.Range("E51").Value = "  -2.59%  "
